$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Gao2022_BrCa
$ws.Range("E4").Value = 44099
$ws.Range("F4").Value = 77.44
$ws.Range("G4").Value = 39913
$ws.Range("H4").Value = 70.09

# Row 5: Gao2022_ERNEG
$ws.Range("E5").Value = 21982
$ws.Range("F5").Value = 77.39
$ws.Range("G5").Value = 19847
$ws.Range("H5").Value = 69.87

# Row 6: Gao2022_ERPOS
$ws.Range("E6").Value = 22685
$ws.Range("F6").Value = 77.4
$ws.Range("G6").Value = 20553
$ws.Range("H6").Value = 70.13
